# Insert a new weekly price-report row for "Ciboulette" (Mercado Mayorista
# Lo Valledor de Santiago) at sheet row 297. This pushes the existing rows
# 297-318 down to 298-319 (dimension grows from R318 to R319), and the new
# row 297 carries the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 297:318 down by inserting a blank row at 297.
$ws.Rows("297:297").Insert()

# Populate the newly-inserted row 297 with the new weekly record.
$ws.Range("A297").Value = 6
$ws.Range("B297").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C297").Value = "Metropolitana"
$ws.Range("D297").Value = 44516
$ws.Range("E297").Value = 13
$ws.Range("F297").Value = 100112039
$ws.Range("G297").Value = "Ciboulette"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 1020
$ws.Range("K297").Value = 700
$ws.Range("L297").Value = 800
$ws.Range("M297").Value = 746
$ws.Range("N297").Value = "$/docena de atados"
$ws.Range("O297").Value = "Región Metropolitana"
$ws.Range("P297").Value = 249
$ws.Range("Q297").Value = 3
$ws.Range("R297").Value = "Hortaliza"
